# "Logic of attendance employee" - correct the attendance rows:
#  - employee id numbers
#  - check-out time values
#  - replace the free-text "date" column with a real date value (no more text)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = 5
$ws.Range("C1").Value = 44542.708333333336
$ws.Range("D1").Value = 44542

# Row 2
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = 44542.625
$ws.Range("D2").Value = 44542

# Move / update the current selection on the sheet
$ws.Range("C1").Select()

# Reposition the workbook window
$win = $excel.ActiveWindow
$win.Top = 1500
$win.Left = 1500
